$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(44313, 2, 6, 113.6363636363636),
    @(44314, 0, 6, 113.6363636363636),
    @(44315, 0, 5, 94.6969696969697),
    @(44316, 1, 5, 94.6969696969697),
    @(44317, 2, 6, 113.6363636363636),
    @(44318, 2, 7, 132.5757575757576)
)

$startRow = 239

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
